$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1928.5
$ws.Range("J18").Value = 2070
$ws.Range("L18").Value = 2070
$ws.Range("N18").Value = -2638
$ws.Range("H33").Value = 257.53333
$ws.Range("I33").Value = 257.53333
$ws.Range("K33").Value = 257.53333
$ws.Range("M33").Value = -28.53332999999998
$ws.Range("H40").Value = 83336600
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 83336600
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 83336600
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -83336950
$ws.Range("H43").Value = 4167.273
$ws.Range("J43").Value = 4658.3335
$ws.Range("L43").Value = 4658.3335
$ws.Range("N43").Value = -4796.3335
$ws.Range("H116").Value = 21285.715
$ws.Range("I116").Value = 11000
$ws.Range("J116").Value = 25400
$ws.Range("K116").Value = 11000
$ws.Range("L116").Value = 25400
$ws.Range("M116").Value = -7558
$ws.Range("N116").Value = -32284
$ws.Range("H135").Value = 1590.2333
$ws.Range("I135").Value = 773.46155
$ws.Range("K135").Value = 6961.15395
$ws.Range("M135").Value = -4426.15395
$ws.Range("H138").Value = 2739.3552
$ws.Range("J138").Value = 4282.0835
$ws.Range("L138").Value = 12846.2505
$ws.Range("N138").Value = -23126.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1676.9166
$ws.Range("I2").Value = 1882.8572
$ws.Range("K2").Value = 1882.8572
$ws.Range("M2").Value = -1769.8572
$ws.Range("H32").Value = 4259.317
$ws.Range("I32").Value = 4482.6
$ws.Range("K32").Value = 4482.6
$ws.Range("M32").Value = -4195.6
$ws.Range("H61").Value = 5682717.5
$ws.Range("I61").Value = 6257145.5
$ws.Range("J61").Value = 2006379.8
$ws.Range("K61").Value = 6257145.5
$ws.Range("L61").Value = 2006379.8
$ws.Range("M61").Value = -6256933.5
$ws.Range("N61").Value = -2006803.8
$ws.Range("H63").Value = 4052.4443
$ws.Range("I63").Value = 4052.4443
$ws.Range("K63").Value = 4052.4443
$ws.Range("M63").Value = -3366.4443
$ws.Range("H66").Value = 4052.4443
$ws.Range("I66").Value = 4052.4443
$ws.Range("K66").Value = 20262.2215
$ws.Range("M66").Value = -16830.2215
$ws.Range("H102").Value = 1889.9166
$ws.Range("I102").Value = 1354.2222
$ws.Range("J102").Value = 3497
$ws.Range("K102").Value = 1354.2222
$ws.Range("L102").Value = 3497
$ws.Range("M102").Value = 267.7778000000001
$ws.Range("N102").Value = -6741
$ws.Range("H116").Value = 1676.9166
$ws.Range("I116").Value = 1882.8572
$ws.Range("K116").Value = 1882.8572
$ws.Range("M116").Value = 411.1428000000001
$ws.Range("H122").Value = 1728.8
$ws.Range("I122").Value = 1244.174
$ws.Range("K122").Value = 3732.522
$ws.Range("M122").Value = -1282.522
$ws.Range("H136").Value = 5682717.5
$ws.Range("I136").Value = 6257145.5
$ws.Range("J136").Value = 2006379.8
$ws.Range("K136").Value = 18771436.5
$ws.Range("L136").Value = 6019139.4
$ws.Range("M136").Value = -18768886.5
$ws.Range("N136").Value = -6024239.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1676.9166
$ws.Range("I3").Value = 1882.8572
$ws.Range("K3").Value = 1882.8572
$ws.Range("M3").Value = -1768.8572
$ws.Range("H20").Value = 2027.4166
$ws.Range("I20").Value = 1564
$ws.Range("K20").Value = 1564
$ws.Range("M20").Value = -1317
$ws.Range("H86").Value = 2604.72
$ws.Range("I86").Value = 2232.8572
$ws.Range("J86").Value = 3078
$ws.Range("K86").Value = 2232.8572
$ws.Range("L86").Value = 3078
$ws.Range("M86").Value = -1109.8572
$ws.Range("N86").Value = -5324
$ws.Range("H89").Value = 2604.72
$ws.Range("I89").Value = 2232.8572
$ws.Range("J89").Value = 3078
$ws.Range("K89").Value = 11164.286
$ws.Range("L89").Value = 15390
$ws.Range("M89").Value = -5548.286
$ws.Range("N89").Value = -26622
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("H105").Value = 1348399.9
$ws.Range("I105").Value = 2858362.5
$ws.Range("J105").Value = 6210.8887
$ws.Range("K105").Value = 2858362.5
$ws.Range("L105").Value = 6210.8887
$ws.Range("M105").Value = -2856615.5
$ws.Range("N105").Value = -9704.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38853612
$ws.Range("I31").Value = 55558930
$ws.Range("K31").Value = 55558930
$ws.Range("M31").Value = -55558635
$ws.Range("H34").Value = 38853612
$ws.Range("I34").Value = 55558930
$ws.Range("K34").Value = 55558930
$ws.Range("M34").Value = -55558728
$ws.Range("H58").Value = 2968.5334
$ws.Range("I58").Value = 3051.75
$ws.Range("J58").Value = 2873.4285
$ws.Range("K58").Value = 3051.75
$ws.Range("L58").Value = 2873.4285
$ws.Range("M58").Value = -2848.75
$ws.Range("N58").Value = -3279.4285
$ws.Range("H68").Value = 85798.664
$ws.Range("J68").Value = 85798.664
$ws.Range("L68").Value = 85798.664
$ws.Range("N68").Value = -87296.664
$ws.Range("H71").Value = 85798.664
$ws.Range("J71").Value = 85798.664
$ws.Range("L71").Value = 257395.992
$ws.Range("N71").Value = -264883.992
$ws.Range("H74").Value = 94443.5
$ws.Range("J74").Value = 94443.5
$ws.Range("L74").Value = 94443.5
$ws.Range("N74").Value = -96191.5
$ws.Range("H77").Value = 94443.5
$ws.Range("J77").Value = 94443.5
$ws.Range("L77").Value = 283330.5
$ws.Range("N77").Value = -292066.5
$ws.Range("H122").Value = 3820.182
$ws.Range("I122").Value = 3660.5715
$ws.Range("J122").Value = 4099.5
$ws.Range("K122").Value = 10981.7145
$ws.Range("L122").Value = 12298.5
$ws.Range("M122").Value = -8531.7145
$ws.Range("N122").Value = -17198.5
$ws.Range("H132").Value = 4309.375
$ws.Range("I132").Value = 2895.2
$ws.Range("J132").Value = 6666.3335
$ws.Range("K132").Value = 8685.599999999999
$ws.Range("L132").Value = 19999.0005
$ws.Range("M132").Value = -6155.599999999999
$ws.Range("N132").Value = -25059.0005
$ws.Range("H134").Value = 4440.6
$ws.Range("I134").Value = 4475.875
$ws.Range("K134").Value = 13427.625
$ws.Range("M134").Value = -10892.625
$ws.Range("H136").Value = 2968.5334
$ws.Range("I136").Value = 3051.75
$ws.Range("J136").Value = 2873.4285
$ws.Range("K136").Value = 9155.25
$ws.Range("L136").Value = 8620.2855
$ws.Range("M136").Value = -6605.25
$ws.Range("N136").Value = -13720.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 16066.6
$ws.Range("J41").Value = 15083.25
$ws.Range("L41").Value = 45249.75
$ws.Range("N41").Value = -45925.75
$ws.Range("H128").Value = 329997.5
$ws.Range("I128").Value = 329997.5
$ws.Range("K128").Value = 989992.5
$ws.Range("M128").Value = -985012.5
$ws.Range("H141").Value = 10362.167
$ws.Range("I141").Value = 5768
$ws.Range("K141").Value = 17304
$ws.Range("M141").Value = -12124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9376.5
$ws.Range("I70").Value = 6506.1665
$ws.Range("J70").Value = 11529.25
$ws.Range("K70").Value = 6506.1665
$ws.Range("L70").Value = 11529.25
$ws.Range("M70").Value = -6236.1665
$ws.Range("N70").Value = -12069.25
$ws.Range("H73").Value = 9376.5
$ws.Range("I73").Value = 6506.1665
$ws.Range("J73").Value = 11529.25
$ws.Range("K73").Value = 6506.1665
$ws.Range("L73").Value = 11529.25
$ws.Range("M73").Value = -5570.1665
$ws.Range("N73").Value = -13401.25
$ws.Range("H102").Value = 3299.2
$ws.Range("I102").Value = 2999.5
$ws.Range("K102").Value = 2999.5
$ws.Range("M102").Value = -1377.5
$ws.Range("H122").Value = 2764.9355
$ws.Range("I122").Value = 2764.9355
$ws.Range("K122").Value = 8294.806500000001
$ws.Range("M122").Value = -5844.806500000001
$ws.Range("H132").Value = 38183196
$ws.Range("I132").Value = 1694.5
$ws.Range("J132").Value = 63637530
$ws.Range("K132").Value = 5083.5
$ws.Range("L132").Value = 190912590
$ws.Range("M132").Value = -2553.5
$ws.Range("N132").Value = -190917650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5179.4116
$ws.Range("I82").Value = 4551.2856
$ws.Range("J82").Value = 5619.1
$ws.Range("K82").Value = 4551.2856
$ws.Range("L82").Value = 5619.1
$ws.Range("M82").Value = -4190.2856
$ws.Range("N82").Value = -6341.1
$ws.Range("H85").Value = 5179.4116
$ws.Range("I85").Value = 4551.2856
$ws.Range("J85").Value = 5619.1
$ws.Range("K85").Value = 4551.2856
$ws.Range("L85").Value = 5619.1
$ws.Range("M85").Value = -3303.2856
$ws.Range("N85").Value = -8115.1
$ws.Range("H136").Value = 4309.5557
$ws.Range("I136").Value = 3763.3914
$ws.Range("K136").Value = 11290.1742
$ws.Range("M136").Value = -8740.174199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 268087.47
$ws.Range("I136").Value = 5094.971
$ws.Range("K136").Value = 15284.913
$ws.Range("M136").Value = -12734.913
